$d = $word.ActiveDocument

function Get-InsertionPoint($para) {
    $ip = $para.Duplicate
    $ip.Start = $para.End - 1
    $ip.End = $para.End - 1
    return $ip
}

function Insert-ColoredText($para, $text, $colorVal) {
    $start = $para.End - 1
    $ip = Get-InsertionPoint $para
    $ip.InsertAfter($text)
    $end = $para.End - 1
    $r = $d.Range($start, $end)
    $r.Font.Color = $colorVal
    return $r
}

# ---- Paragraph 1: "This is a Microsoft word document." ----
$p1 = $d.Paragraphs(1).Range

# Run 1 (plain): append two trailing spaces
$ip1 = Get-InsertionPoint $p1
$ip1.InsertAfter("  ")

$RED = 192  # 0xC00000 packed as BGR -> 192

$enDash = [char]0x2013
$text2 = "(This is a change " + $enDash + " Ve"
$text3 = "rsion for branch alternate"
$text4 = ")"

Insert-ColoredText $p1 $text2 $RED | Out-Null
Insert-ColoredText $p1 $text3 $RED | Out-Null
Insert-ColoredText $p1 $text4 $RED | Out-Null

Write-Host ("P1 final text: [" + $p1.Text + "]")

# ---- After paragraph 2, insert a new empty paragraph with specific formatting ----
$p2 = $d.Paragraphs(2).Range
$newPara = $p2.InsertParagraphAfter()
